$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price column (D): force Text format so numeric-looking strings
#     (e.g. "0.08400", "6.090") keep their exact digits/trailing zeros
#     instead of being auto-converted to numbers. ---
$priceCells = @("D2","D3","D4","D5","D7","D8","D9","D13","D14","D15","D16","D17","D19","D20","D22","D23","D24","D25","D26","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.246.91'
$ws.Range("D3").Value = '1.895.41'
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").Value = '325.05'
$ws.Range("D7").Value = '0.5178'
$ws.Range("D8").Value = '0.4008'
$ws.Range("D9").Value = '0.08400'
$ws.Range("D13").Value = '6.438'
$ws.Range("D14").Value = '1.893.34'
$ws.Range("D15").Value = '7.322'
$ws.Range("D16").Value = '1.002'
$ws.Range("D17").Value = '94.16'
$ws.Range("D19").Value = '0.06641'
$ws.Range("D20").Value = '18.21'
$ws.Range("D22").Value = '5.951'
$ws.Range("D23").Value = '30.232.47'
$ws.Range("D24").Value = '11.29'
$ws.Range("D25").Value = '2.229'
$ws.Range("D26").Value = '2.110.70'
$ws.Range("D28").Value = '161.77'
$ws.Range("D29").Value = '2.353'
$ws.Range("D30").Value = '129.52'
$ws.Range("D32").Value = '0.1055'
$ws.Range("D33").Value = '6.090'
$ws.Range("D34").Value = '3.744'
$ws.Range("D35").Value = '0.02494'
$ws.Range("D36").Value = '0.06548'
$ws.Range("D37").Value = '5.284'
$ws.Range("D38").Value = '0.2196'
$ws.Range("D39").Value = '1.223'
$ws.Range("D41").Value = '8.741'
$ws.Range("D42").Value = '0.6498'
$ws.Range("D43").Value = '1.228'
$ws.Range("D44").Value = '0.6097'
$ws.Range("D45").Value = '13.21'
$ws.Range("D46").Value = '3.696'
$ws.Range("D47").Value = '2.054'
$ws.Range("D48").Value = '1.235'
$ws.Range("D49").Value = '124.52'

foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "General"
}

# --- Update Volume(1h) column (E) ---
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +3.10%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("E12").Value = '  +10.70%  '
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("E27").Value = '  +2.96%  '
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("E29").Value = '  -3.70%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("E34").Value = '  +2.73%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("E41").Value = '  -3.96%  '
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("E51").Value = '  +1.19%  '
